$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Project IDs that should NOT receive a Flag_path value
$excludedIds = @(790, 1606, 1700, 4023, 6111, 6729, 3274)

for ($row = 2; $row -le 95; $row++) {
    $projectId = $ws.Cells.Item($row, 1).Value2
    if ($null -eq $projectId -or $projectId -eq "") { continue }
    if ($excludedIds -contains $projectId) { continue }
    $path = "~/GIT/PSSdb/raw/flags/Zooscan/project_{0}_flags.tsv" -f $projectId
    $ws.Cells.Item($row, 37).Value2 = $path
}
